$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Header text tweaks (case changes)
$ws.Range("I1").Value = "TOTAL PHYSICAL TARGET"
$ws.Range("L1").Value = "BATCH"

# 2) Insert 5 new columns before the existing "Status as of July 4, 2025"
#    column (currently AA). This shifts that column (and its data
#    validation) to AF automatically, and the newly inserted blank
#    columns inherit the header styling from their left neighbor (Z1).
$ws.Range("AA1:AE1").EntireColumn.Insert()

# 3) Fill in the headers for the 5 newly inserted columns
$ws.Range("AA1").Value = "No. of Sites Reverted"
$ws.Range("AB1").Value = "No. of Sites Not yet started"
$ws.Range("AC1").Value = "No. of Sites Under Procurement"
$ws.Range("AD1").Value = "No. of Sites On Going"
$ws.Range("AE1").Value = "No. of Sites Completed"

# 4) Remove the placeholder "-" entries left in columns I (Total Physical
#    Target) and L (Batch) for the first 21 data rows (rows 2-22).
for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 9).ClearContents()
    $ws.Cells.Item($r, 12).ClearContents()
}
